# Daily attendance processing - normalize "Recorded By" (column G) lists.
#
# Each "Recorded By" cell holds a comma-separated list of recorder names /
# emails (e.g. "System, dnasr281@gmail.com"). This pass re-sorts each list
# alphabetically (case-insensitive), breaking ties between entries that only
# differ by case so that the lower-case spelling sorts before the
# upper-case one (e.g. "system" before "System").

function Compare-Recorder($a, $b) {
    $la = $a.ToLower()
    $lb = $b.ToLower()
    if ($la -lt $lb) { return -1 }
    if ($la -gt $lb) { return 1 }
    # Case-insensitive tie: fall back to a case-sensitive comparison, but
    # invert it so the lower-case variant sorts first.
    $c = $a.CompareTo($b)
    if ($c -eq 0) { return 0 }
    if ($c -lt 0) { return 1 }
    return -1
}

function Sort-Recorders($arr) {
    $n = $arr.Count
    for ($i = 0; $i -lt $n; $i++) {
        for ($j = 0; $j -lt ($n - $i - 1); $j++) {
            $cmp = Compare-Recorder $arr[$j] $arr[$j + 1]
            if ($cmp -gt 0) {
                $tmp = $arr[$j]
                $arr[$j] = $arr[$j + 1]
                $arr[$j + 1] = $tmp
            }
        }
    }
    return $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

$recordedByCol = 7   # column G = "Recorded By"

for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $val = $cell.Text

    if ($val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $sortedParts = Sort-Recorders $parts
            $newVal = $sortedParts -join ", "
            if (-not $newVal.Equals($val)) {
                $cell.Value = $newVal
            }
        }
    }
}
